Write-Host "noop"
